{"js": "// Update the date line and the 25 division problems in the practice table.\n// Each (row, col) table-cell coordinate is addressed directly so that\n// duplicate problem text (e.g. \"17\u00f75=\" appears both as an old value in one\n// cell and a new value in another) never causes ambiguous matches.\n\n// 1) Date heading paragraph: \"2024-06-15 Saturday\" -> \"2024-06-16 Sunday\"\nconst body = context.document.body;\nconst dateResults = body.search(\"2024-06-15 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2024-06-16 Sunday\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Division problems laid out in a 20-row x 5-col table; only rows\n//    0, 4, 8, 12, 16 hold text (the rest are blank spacer rows).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellEdits = [\n  [0, 0, \"48\u00f74=\", \"59\u00f75=\"],\n  [0, 1, \"14\u00f73=\", \"37\u00f72=\"],\n  [0, 2, \"90\u00f79=\", \"82\u00f75=\"],\n  [0, 3, \"27\u00f72=\", \"19\u00f78=\"],\n  [0, 4, \"35\u00f72=\", \"20\u00f74=\"],\n  [4, 0, \"79\u00f72=\", \"64\u00f72=\"],\n  [4, 1, \"93\u00f78=\", \"48\u00f75=\"],\n  [4, 2, \"87\u00f72=\", \"32\u00f73=\"],\n  [4, 3, \"17\u00f74=\", \"92\u00f79=\"],\n  [4, 4, \"66\u00f72=\", \"82\u00f78=\"],\n  [8, 0, \"69\u00f72=\", \"17\u00f75=\"],\n  [8, 1, \"65\u00f73=\", \"51\u00f73=\"],\n  [8, 2, \"83\u00f75=\", \"18\u00f78=\"],\n  [8, 3, \"78\u00f72=\", \"95\u00f73=\"],\n  [8, 4, \"14\u00f74=\", \"34\u00f78=\"],\n  [12, 0, \"37\u00f79=\", \"23\u00f78=\"],\n  [12, 1, \"17\u00f75=\", \"65\u00f79=\"],\n  [12, 2, \"59\u00f78=\", \"46\u00f72=\"],\n  [12, 3, \"73\u00f74=\", \"92\u00f72=\"],\n  [12, 4, \"65\u00f76=\", \"94\u00f76=\"],\n  [16, 0, \"33\u00f72=\", \"13\u00f74=\"],\n  [16, 1, \"22\u00f78=\", \"20\u00f79=\"],\n  [16, 2, \"82\u00f74=\", \"10\u00f72=\"],\n  [16, 3, \"98\u00f78=\", \"29\u00f79=\"],\n  [16, 4, \"24\u00f78=\", \"95\u00f76=\"],\n];\n\nfor (const [row, col, oldText, newText] of cellEdits) {\n  const cell = table.getCell(row, col);\n  const cellBody = cell.body;\n  const results = cellBody.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice table.\n# Each edit targets an exact (row, col) table-cell coordinate (1-based, as\n# COM indexes tables/cells) instead of a document-wide Find/Replace so that\n# duplicate problem text (e.g. \"17\u00f75=\" is an old value in one cell and the\n# new value written into a different cell) never causes an unintended\n# second replacement.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph: \"2024-06-15 Saturday\" -> \"2024-06-16 Sunday\"\n$d.Paragraphs.Item(1).Range.Text = \"2024-06-16 Sunday\"\n\n# 2) Division problems laid out in a 20-row x 5-col table; only rows\n#    1, 5, 9, 13, 17 (1-based) hold text - the rest are blank spacer rows.\n$tbl = $d.Tables.Item(1)\n\n$edits = @(\n  @{ Row=1;  Col=1; Old=\"48\u00f74=\"; New=\"59\u00f75=\" },\n  @{ Row=1;  Col=2; Old=\"14\u00f73=\"; New=\"37\u00f72=\" },\n  @{ Row=1;  Col=3; Old=\"90\u00f79=\"; New=\"82\u00f75=\" },\n  @{ Row=1;  Col=4; Old=\"27\u00f72=\"; New=\"19\u00f78=\" },\n  @{ Row=1;  Col=5; Old=\"35\u00f72=\"; New=\"20\u00f74=\" },\n  @{ Row=5;  Col=1; Old=\"79\u00f72=\"; New=\"64\u00f72=\" },\n  @{ Row=5;  Col=2; Old=\"93\u00f78=\"; New=\"48\u00f75=\" },\n  @{ Row=5;  Col=3; Old=\"87\u00f72=\"; New=\"32\u00f73=\" },\n  @{ Row=5;  Col=4; Old=\"17\u00f74=\"; New=\"92\u00f79=\" },\n  @{ Row=5;  Col=5; Old=\"66\u00f72=\"; New=\"82\u00f78=\" },\n  @{ Row=9;  Col=1; Old=\"69\u00f72=\"; New=\"17\u00f75=\" },\n  @{ Row=9;  Col=2; Old=\"65\u00f73=\"; New=\"51\u00f73=\" },\n  @{ Row=9;  Col=3; Old=\"83\u00f75=\"; New=\"18\u00f78=\" },\n  @{ Row=9;  Col=4; Old=\"78\u00f72=\"; New=\"95\u00f73=\" },\n  @{ Row=9;  Col=5; Old=\"14\u00f74=\"; New=\"34\u00f78=\" },\n  @{ Row=13; Col=1; Old=\"37\u00f79=\"; New=\"23\u00f78=\" },\n  @{ Row=13; Col=2; Old=\"17\u00f75=\"; New=\"65\u00f79=\" },\n  @{ Row=13; Col=3; Old=\"59\u00f78=\"; New=\"46\u00f72=\" },\n  @{ Row=13; Col=4; Old=\"73\u00f74=\"; New=\"92\u00f72=\" },\n  @{ Row=13; Col=5; Old=\"65\u00f76=\"; New=\"94\u00f76=\" },\n  @{ Row=17; Col=1; Old=\"33\u00f72=\"; New=\"13\u00f74=\" },\n  @{ Row=17; Col=2; Old=\"22\u00f78=\"; New=\"20\u00f79=\" },\n  @{ Row=17; Col=3; Old=\"82\u00f74=\"; New=\"10\u00f72=\" },\n  @{ Row=17; Col=4; Old=\"98\u00f78=\"; New=\"29\u00f79=\" },\n  @{ Row=17; Col=5; Old=\"24\u00f78=\"; New=\"95\u00f76=\" }\n)\n\nforeach ($edit in $edits) {\n  $cell = $tbl.Cell($edit.Row, $edit.Col)\n  $cell.Range.Text = $edit.New\n}\n\nWrite-Output \"done\"\n"}
